# Updated model; restructured data mapping
# Target sheet: "attributes" (the 3rd worksheet) holds attribute metadata
# for the rd3stats_treedata entity. We add a new "id" attribute row and
# rework the nillable/labelAttribute/lookupAttribute/dataType columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attributes")

# --- Header row: columns E-H were restructured. ---
# Old order: idAttribute, labelAttribute, lookupAttribute, nillable, dataType
# New order: idAttribute, nillable, dataType, labelAttribute, lookupAttribute
$ws.Range("E1").Value = "nillable"
$ws.Range("F1").Value = "dataType"
$ws.Range("G1").Value = "labelAttribute"
$ws.Range("H1").Value = "lookupAttribute"

# --- Row 2: id attribute (brand new semantics, re-using row 2) ---
$ws.Range("A2").Value = "rd3stats_treedata"
$ws.Range("B2").Value = "id"
$ws.Range("C2").Value = "identifier of the object"
$ws.Range("D2").Value = $true
$ws.Range("E2").Value = $false
$ws.Range("F2").Value = "string"
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false

# --- Row 3: subjectID ---
$ws.Range("A3").Value = "rd3stats_treedata"
$ws.Range("B3").Value = "subjectID"
$ws.Range("C3").Value = "An individual who is the subject of personal data, persons to whom data refers, and from whom data are collected, processed, and stored."
$ws.Range("D3").Value = $false
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = "string"
$ws.Range("G3").Value = $true
$ws.Range("H3").Value = $true

# --- Row 4: familyID ---
$ws.Range("A4").Value = "rd3stats_treedata"
$ws.Range("B4").Value = "familyID"
$ws.Range("C4").Value = "A domestic group, or a number of domestic groups linked through descent (demonstrated or stipulated) from a common ancestor, marriage, or adoption."
$ws.Range("D4").Value = $false
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = "string"
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = $false

# --- Row 5: json (new row) ---
$ws.Range("A5").Value = "rd3stats_treedata"
$ws.Range("B5").Value = "json"
$ws.Range("C5").Value = "json stringified object containing sample-experiment links"
$ws.Range("D5").Value = $false
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = "text"
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = $false
